$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 19712
$ws1.Range("F5").Value = 0
$ws1.Range("F7").Value = 0
$ws1.Range("F9").Value = 0
$ws1.Range("F12").Value = 0
$ws1.Range("F13").Value = 0
$ws1.Range("F15").Value = 0
$ws1.Range("F16").Value = 0
$ws1.Range("F18").Value = 188
$ws1.Range("F21").Value = 0
$ws1.Range("F22").Value = 0
$ws1.Range("F27").Value = 1081
$ws1.Range("F30").Value = 0
$ws1.Range("F32").Value = 0
$ws1.Range("F33").Value = 51
$ws1.Range("F34").Value = 2790
$ws1.Range("F35").Value = 0
$ws1.Range("F37").Value = 19
$ws1.Range("F38").Value = 12549
$ws1.Range("F40").Value = 64
$ws1.Range("F41").Value = 0
$ws1.Range("F43").Value = 0
$ws1.Range("F44").Value = 0
$ws1.Range("F45").Value = 0

# Sheet "全部类型" (sheet4.xml)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 0
$ws4.Range("F3").Value = 0
$ws4.Range("F4").Value = 19712
$ws4.Range("F5").Value = 794
$ws4.Range("F7").Value = 0
$ws4.Range("F8").Value = 0
$ws4.Range("F9").Value = 7458
$ws4.Range("F13").Value = 0
$ws4.Range("F14").Value = 0
$ws4.Range("F16").Value = 0
$ws4.Range("F18").Value = 0
$ws4.Range("F21").Value = 0
$ws4.Range("F22").Value = 0
$ws4.Range("F24").Value = 0
$ws4.Range("F26").Value = 0
$ws4.Range("F27").Value = 0
$ws4.Range("F28").Value = 0
$ws4.Range("F30").Value = 172
$ws4.Range("F31").Value = 0
$ws4.Range("F32").Value = 0
$ws4.Range("F34").Value = 51
$ws4.Range("F36").Value = 0
$ws4.Range("F38").Value = 0
$ws4.Range("F39").Value = 0
$ws4.Range("F40").Value = 12549
$ws4.Range("F42").Value = 64
$ws4.Range("F43").Value = 0
$ws4.Range("F44").Value = 0
$ws4.Range("F45").Value = 0
$ws4.Range("F46").Value = 350
